# Updated symbol list on Fri Jan 20 19:38:58 UTC 2023 with GitHub Actions
#
# Refresh the Price (col D) and Volume(1h) (col E) columns for the rows
# whose quote changed. Values are stored as plain text in the sheet
# (e.g. "292.75", "-0.01%"), so each cell is briefly switched to the
# Text number format before the assignment (otherwise Excel would parse
# a numeric-looking string, or a "N%" string, into a real number/percent
# and reformat the cell) and then ClearFormats() restores the cell's
# original (default) style so only the cell's value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# row -> @(newPrice, newVolume)   ($null entries are left untouched)
$changes = @{
    2  = @("292.75",     "-0.01%")
    3  = @("31.23",      "1.22%")
    4  = @("4.980",      "1.60%")
    5  = @("0.07497",    "2.69%")
    6  = @("2.282",      "-1.80%")
    7  = @("7.790",      "1.27%")
    8  = @("0.9204",     "2.25%")
    9  = @("0.09267",    "16.94%")
    10 = @("0.1730",     "2.78%")
    11 = @("0.08348",    "2.10%")
    12 = @("0.03282",    "6.05%")
    13 = @("0.09948",    "-1.01%")
    14 = @("0.001503",   "-1.39%")
    15 = @("0.005784",   "0.24%")
    16 = @($null,        "0.13%")
    17 = @("3.772",      "1.13%")
    18 = @("2.166",      "4.35%")
    19 = @("0.3343",     "0.42%")
    20 = @("0.1315",     "1.13%")
    21 = @("4.105",      "1.77%")
    22 = @("0.2101",     "0.27%")
    23 = @("0.04540",    "0.03%")
    24 = @("0.001220",   "0.92%")
    25 = @("0.004306",   "-7.02%")
    26 = @("0.0001299",  "0.18%")
    27 = @("0.0003396",  "0.29%")
    39 = @($null,        "1.84%")
    40 = @("0.04590",    "3.95%")
    41 = @("0.007459",   "2.28%")
    42 = @("0.009845",   "14.49%")
    43 = @("0.1359",     "2.94%")
    44 = @("0.002158",   "8.07%")
    45 = @("0.009780",   "3.92%")
    46 = @("0.00006113", "3.27%")
    47 = @($null,        "0.18%")
    48 = @("2.654",      "18.45%")
    49 = @("0.002001",   "-30.79%")
    50 = @("0.00002099", "0.18%")
    51 = @("0.0001999",  "0.18%")
}

foreach ($r in $changes.Keys) {
    $pair = $changes[$r]
    $newPrice  = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        Set-TextValue $ws $r 4 $newPrice
    }
    if ($null -ne $newVolume) {
        Set-TextValue $ws $r 5 $newVolume
    }
}
